$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell K1 (matches style of existing header row, e.g. J1)
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("K1").Value = "PhylogenySorting"

# New data cell K4
$ws.Range("K4").Value = "T45"

# Column widths for J (10) and K (11) -- closest achievable to 12.6640625 / 17.6640625
$ws.Columns.Item(10).ColumnWidth = 11.83
$ws.Columns.Item(11).ColumnWidth = 16.83

# Update selection to match diff
$ws.Range("G8").Select() | Out-Null
